$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.769.55'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.089.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.19%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.78'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.66'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.93%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0838'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.92%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.397.59'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.98'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.91'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.797'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.15%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.095.01'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.98%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.699.23'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.64'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.03'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0838'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.82'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.39%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.44%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.05'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.97%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.46'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.138'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.45'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +12.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.17'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.121'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.35'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.51'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.87%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.77%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.44'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.21%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.55'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.17%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.21'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.542.70'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.82%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.98'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.14%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0225'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.31%  '

$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0927'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.19%  '

$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.82'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.73%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +8.57%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.10'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.60%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.40%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.285.88'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.26%  '
